$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.314.49"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -7.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.493.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.73%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "394.84"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.75"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.483.07"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -8.96%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.691"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -11.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -14.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000291"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -14.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.95"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -6.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.43"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.065.65"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.478.00"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.12"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.73"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.58%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -8.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "62.259.68"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -7.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "404.92"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -10.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.23"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +8.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.89"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -8.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.87"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -9.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "33.98"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.06"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -8.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.12"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.05%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -10.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.05"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.97%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.111"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.49%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.62"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -10.52%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.154"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "38.87"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.69"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0444"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -10.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.994"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.75"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +19.33%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.135"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -8.75%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.95"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.33%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.00"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +16.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.15"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.73%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.98"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.14"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.22%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.15"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +16.80%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.71"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -11.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.48"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -9.69%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "PEPE"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0588"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -21.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.276"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -12.90%  "
